$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVERYTHING")
$ws.Range("A2").Interior.Color = 16777215
$ws.Range("A2").Font.Bold = $true
